# Sync attendance_reports: add new Attendance check-ins (rows 395-414) and
# refresh the computed Summary rows for the students affected by them.

$wb = $excel.ActiveWorkbook
$wsAtt = $wb.Worksheets.Item("Attendance")
$wsSum = $wb.Worksheets.Item("Summary")

# ---------------------------------------------------------------------------
# 1) Attendance sheet: append 20 new ANATOMY check-in rows (395-414)
# ---------------------------------------------------------------------------

$newRows = @(
    @("221579","لينا مكرم محمد يسن","Year 2","C1","221579@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:23:11","C1"),
    @("221574","عمر براء رجب","Year 2","C1","221574@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:23:17","C1"),
    @("221654","محمد اسامه بابكر احمد","Year 2","C1","221654@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:23:22","C1"),
    @("221566","مصطفى سامى محمد عبد الله","Year 2","C1","221566@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:23:27","C1"),
    @("221506","فاطمه احمد اسماعيل الناجي","Year 2","C1","221506@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:23:39","C1"),
    @("220990","جاد زياد سلوم","Year 2","C1","220990@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:23:48","C1"),
    @("222004","احمد ايمن احمد بشير","Year 2","C1","222004@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:23:52","C1"),
    @("210935","يعقوب سليمان يعقوب يحى","Year 2","C1","210935@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:23:56","C1"),
    @("221653","يزن يحيى سليمان طبش","Year 2","C1","221653@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:24:10","C1"),
    @("221546","محمدزين ابوبكر محمد زين احمد","Year 2","C1","221546@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:24:20","C1"),
    @("221863","ميسم ايمن زيدان","Year 2","C1","221863@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:24:31","C1"),
    @("221802","بشاير ابوبكر على عيسى ابراهيم","Year 2","C1","221802@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:24:42","C1"),
    @("210967","ملاك كمال اسماعيل ابو جلاله","Year 2","C1","210967@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:24:57","C1"),
    @("221357","عبد الله محمد نصر قناوى","Year 2","C1","221357@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:38:41","C1"),
    @("211915","ابايزيد عبد الله سعيد ابو رصاص","Year 2","C1","211915@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:38:49","C1"),
    @("221031","امنيه عبدالله عبد اللطيف محمد","Year 2","C1","221031@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:40:56","C1"),
    @("212318","مازن مصدق يس عبد اللطيف","Year 2","C1","212318@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:41:03","C1"),
    @("221433","عادل سامي احمد طه","Year 2","C1","221433@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:42:09","C1"),
    @("221437","مهند عدنان دخل الله ماضي","Year 2","C1","221437@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:43:37","C1"),
    @("221997","خالد مبارك احمد عثمان","Year 2","C1","221997@med.asu.edu.eg","ANATOMY","2","ANATOMY","09/11/2025","14:46:28","C1")
)

$startRow = 395
$endRow = $startRow + $newRows.Count - 1

# Pre-format the whole new block as Text so numeric-looking values (student
# IDs, the "2" session number, the dd/mm/yyyy date) are stored as literal
# strings instead of being auto-coerced into numbers / date serials.
$wsAtt.Range("A$startRow`:K$endRow").NumberFormat = "@"

$r = $startRow
foreach ($row in $newRows) {
    $wsAtt.Cells.Item($r, 1).Value = $row[0]
    $wsAtt.Cells.Item($r, 2).Value = $row[1]
    $wsAtt.Cells.Item($r, 3).Value = $row[2]
    $wsAtt.Cells.Item($r, 4).Value = $row[3]
    $wsAtt.Cells.Item($r, 5).Value = $row[4]
    $wsAtt.Cells.Item($r, 6).Value = $row[5]
    $wsAtt.Cells.Item($r, 7).Value = $row[6]
    $wsAtt.Cells.Item($r, 8).Value = $row[7]
    $wsAtt.Cells.Item($r, 9).Value = $row[8]
    $wsAtt.Cells.Item($r, 10).Value = $row[9]
    $wsAtt.Cells.Item($r, 11).Value = $row[10]
    $r = $r + 1
}

# Keep the sheet's AutoFilter live over the grown range.
$wsAtt.AutoFilterMode = $false
$wsAtt.Range("A1:K$endRow").AutoFilter()

# Extend the hidden _FilterDatabase defined name for the Attendance sheet
# (mirrors Excel's own behaviour when AutoFilter range grows).
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -eq "Attendance!_FilterDatabase") {
        $nm.RefersTo = "=Attendance!`$A`$1:`$K`$$endRow"
    }
}

# ---------------------------------------------------------------------------
# 2) Summary sheet: refresh the derived stats for the 20 affected students
# ---------------------------------------------------------------------------

# Fill colors (OLE BGR ints) matching the workbook's existing Status styles.
$colorLowRisk = 10940927   # FFF1A6
$colorNoRisk = 13994035    # 3388D5
$colorModerateRisk = 8239615   # FFB97D
$colorHighRisk = 8158463    # FF7C7C

# row, newStatus ("" = unchanged), G%, H, L, M, O, Q
$summaryUpdates = @(
    @(22, "No Risk", "17.2%", 17, 5, 1, 2, 1),
    @(23, "", "6.9%", 20, 2, 4, 1, 1),
    @(44, "", "13.8%", 18, 4, 2, 2, 1),
    @(63, "", "13.8%", 18, 4, 2, 2, 1),
    @(95, "", "13.8%", 18, 4, 2, 2, 1),
    @(100, "", "13.8%", 18, 4, 2, 1, 1),
    @(115, "", "13.8%", 18, 4, 2, 2, 1),
    @(125, "", "6.9%", 20, 2, 4, 2, 1),
    @(127, "", "6.9%", 20, 2, 4, 2, 1),
    @(138, "Low Risk", "10.3%", 19, 3, 3, 2, 1),
    @(148, "Low Risk", "10.3%", 19, 3, 3, 1, 1),
    @(151, "Moderate Risk", "3.4%", 21, 1, 5, 1, 1),
    @(153, "", "13.8%", 18, 4, 2, 1, 1),
    @(155, "", "13.8%", 18, 4, 2, 2, 1),
    @(171, "Low Risk", "10.3%", 19, 3, 3, 1, 1),
    @(172, "", "6.9%", 20, 2, 4, 1, 1),
    @(191, "Low Risk", "10.3%", 19, 3, 3, 1, 1),
    @(201, "Low Risk", "10.3%", 19, 3, 3, 1, 1),
    @(229, "", "6.9%", 20, 2, 4, 1, 1),
    @(233, "No Risk", "17.2%", 17, 5, 1, 2, 1)
)

foreach ($u in $summaryUpdates) {
    $row = $u[0]
    $status = $u[1]
    $pct = $u[2]
    $hVal = $u[3]
    $lVal = $u[4]
    $mVal = $u[5]
    $oVal = $u[6]
    $qVal = $u[7]

    if ($status -ne "") {
        $fCell = $wsSum.Cells.Item($row, 6)
        $fCell.Value = $status
        if ($status -eq "No Risk") {
            $fCell.Interior.Color = $colorNoRisk
        } elseif ($status -eq "Low Risk") {
            $fCell.Interior.Color = $colorLowRisk
        } elseif ($status -eq "Moderate Risk") {
            $fCell.Interior.Color = $colorModerateRisk
        } elseif ($status -eq "High Risk") {
            $fCell.Interior.Color = $colorHighRisk
        }
    }

    # The Percentage column stores its text (e.g. "17.2%") literally rather
    # than as a numeric percentage, so force Text format before assigning.
    $gCell = $wsSum.Cells.Item($row, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $pct
    $wsSum.Cells.Item($row, 8).Value = $hVal
    $wsSum.Cells.Item($row, 12).Value = $lVal
    $wsSum.Cells.Item($row, 13).Value = $mVal
    $wsSum.Cells.Item($row, 15).Value = $oVal
    $wsSum.Cells.Item($row, 17).Value = $qVal
}
